$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.595.50'
$ws.Range('E2').Value = '  -3.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.850.58'
$ws.Range('E3').Value = '  -3.94%  '
$ws.Range('E4').Value = '  -1.04%  '
$ws.Range('E5').Value = '  +3.16%  '
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4662'
$ws.Range('E7').Value = '  -3.23%  '
$ws.Range('E8').Value = '  -3.60%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.19'
$ws.Range('E9').Value = '  -3.14%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07905'
$ws.Range('E10').Value = '  -3.65%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9824'
$ws.Range('E11').Value = '  -2.77%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.33'
$ws.Range('E12').Value = '  -6.21%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.878.53'
$ws.Range('E13').Value = '  -2.97%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.847'
$ws.Range('E14').Value = '  -4.17%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.005'
$ws.Range('E15').Value = '  -4.12%  '
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06919'
$ws.Range('E16').Value = '  +0.44%  '
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.003'
$ws.Range('E17').Value = '  -1.05%  '
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '87.67'
$ws.Range('E18').Value = '  -4.24%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001004'
$ws.Range('E19').Value = '  -3.42%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.13'
$ws.Range('E20').Value = '  -3.10%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  -0.96%  '
$ws.Range('B22').Value = 'WrappedBTC'
$ws.Range('C22').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '28.601.91'
$ws.Range('E22').Value = '  -3.25%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.403'
$ws.Range('E23').Value = '  -4.79%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.33'
$ws.Range('E24').Value = '  -5.79%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.162'
$ws.Range('E25').Value = '  -0.62%  '
$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.085.11'
$ws.Range('E26').Value = '  -3.79%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '153.98'
$ws.Range('E27').Value = '  -1.21%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.47'
$ws.Range('E28').Value = '  -2.86%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.120'
$ws.Range('E29').Value = '  -4.53%  '
$ws.Range('B30').Value = 'LidoDAOToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.030'
$ws.Range('E30').Value = '  -3.08%  '
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '117.69'
$ws.Range('E31').Value = '  -2.45%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9745'
$ws.Range('E32').Value = '  -3.91%  '
$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09382'
$ws.Range('E33').Value = '  -2.11%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.366'
$ws.Range('E34').Value = '  -4.08%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.486'
$ws.Range('E35').Value = '  -2.10%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.347'
$ws.Range('E36').Value = '  -2.62%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06160'
$ws.Range('E37').Value = '  -3.09%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02202'
$ws.Range('E38').Value = '  -3.75%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.162'
$ws.Range('E39').Value = '  -2.42%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5729'
$ws.Range('E40').Value = '  -3.77%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.696'
$ws.Range('E41').Value = '  -2.51%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '10.19'
$ws.Range('E42').Value = '  -4.85%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1801'
$ws.Range('E43').Value = '  -2.35%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.372'
$ws.Range('E44').Value = '  -4.24%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.248'
$ws.Range('E45').Value = '  -2.50%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5401'
$ws.Range('E46').Value = '  -2.73%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '11.82'
$ws.Range('E47').Value = '  -5.23%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.07127'
$ws.Range('E48').Value = '  -4.72%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.910'
$ws.Range('E49').Value = '  -3.38%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '115.09'
$ws.Range('E50').Value = '  -2.88%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.06'
$ws.Range('E51').Value = '  +1.54%  '
